$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host ("Sheet name: " + $ws.Name)
Write-Host ("Cell A1: " + $ws.Range("A1").Value)
